$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Prefab (column D) paths for each building row
$ws.Range("D2").Value = "Prefabs/Object/Altar_1_1"
$ws.Range("D3").Value = "Prefabs/Object/Arena_1_1"
$ws.Range("D4").Value = "Prefabs/Object/Camp_1_1"
$ws.Range("D5").Value = "Prefabs/Object/GoldMine_1_1"
$ws.Range("D6").Value = "Prefabs/Object/Item_hourse_1_1"
$ws.Range("D7").Value = "Prefabs/Object/League_1_1"
$ws.Range("D8").Value = "Prefabs/Object/MagicHourse_1_1"
$ws.Range("D9").Value = "Prefabs/Object/Tower_1_1"
$ws.Range("D10").Value = "Prefabs/Object/Town_1_1"

# Move the active selection
$ws.Range("D11").Select()
